$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the "LastName"/"FirstName" header cells (B1 <-> C1)
$b1 = $ws.Range("B1").Value2
$c1 = $ws.Range("C1").Value2
$ws.Range("B1").Value = $c1
$ws.Range("C1").Value = $b1

# Re-apply the default ("Normal") style explicitly on the CNE/LastName/FirstName
# columns for the header row and all data rows
$ws.Range("A1:C11").Style = "Normal"

# Update the CNE identifiers in column A (rows 2-11): 190000xx -> 170000xx
for ($r = 2; $r -le 11; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value = $cell.Value2 - 2000000
}

# Update the last selected cell in the sheet view
$ws.Range("G7").Select()
